$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string texts (order doesn't matter for COM, Excel will manage
# the shared-string table itself based on usage order of new distinct text)
$title = "#dig Transform test. "
$dChar = "d"
$hashChar = "#"

# Clear out the old sheet contents first (old layout used columns A,C:I / rows 1:10)
$ws.Cells.Clear() | Out-Null

# Row 1: just the title, now trimmed (moved conceptually but stays at A1)
$ws.Range("A1").Value = $title

# Row 2: d fills A2:I2, border char in J2
$ws.Range("A2:I2").Value = $dChar
$ws.Range("J2").Value = $hashChar

# Row 3: A3 = d, I3 = d, J3 = #
$ws.Range("A3").Value = $dChar
$ws.Range("I3").Value = $dChar
$ws.Range("J3").Value = $hashChar

# Row 4: A4 = d, I4 = d, J4 = #
$ws.Range("A4").Value = $dChar
$ws.Range("I4").Value = $dChar
$ws.Range("J4").Value = $hashChar

# Row 5: A5 = d, I5 = d, J5 = #
$ws.Range("A5").Value = $dChar
$ws.Range("I5").Value = $dChar
$ws.Range("J5").Value = $hashChar

# Row 6: A6 = d, I6 = d, J6 = #
$ws.Range("A6").Value = $dChar
$ws.Range("I6").Value = $dChar
$ws.Range("J6").Value = $hashChar

# Row 7: A7 = d, E7 = d, I7 = d, J7 = #
$ws.Range("A7").Value = $dChar
$ws.Range("E7").Value = $dChar
$ws.Range("I7").Value = $dChar
$ws.Range("J7").Value = $hashChar

# Row 8: A8 = d, D8:F8 = d, I8 = d, J8 = #
$ws.Range("A8").Value = $dChar
$ws.Range("D8:F8").Value = $dChar
$ws.Range("I8").Value = $dChar
$ws.Range("J8").Value = $hashChar

# Row 9: A9 = d, C9:G9 = d, I9 = d, J9 = #
$ws.Range("A9").Value = $dChar
$ws.Range("C9:G9").Value = $dChar
$ws.Range("I9").Value = $dChar
$ws.Range("J9").Value = $hashChar

# Row 10: A10:I10 = d, J10 = #
$ws.Range("A10:I10").Value = $dChar
$ws.Range("J10").Value = $hashChar

# Row 11: A11:J11 = # (new border row)
$ws.Range("A11:J11").Value = $hashChar
$ws.Rows.Item(11).RowHeight = 17.25

$ws.Range("A2").Select() | Out-Null
